$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# --- Row 1 ---
$ws.Range("N1").Formula = "=INDEX(OFFSET(Data!A1, 0, 0, 3, 3), 2, 2)"
# O1 formula is unchanged: =SUM(INDEX(Data!A1:E6, 0, 2))
$ws.Range("P1").Value = "Data!B2"
$ws.Range("Q1").Value = 25
$ws.Range("Z1").Value = "Test Value"

# --- Row 2 ---
$ws.Range("N2").Formula = "=OFFSET(INDEX(Data!A1:E6, 2, 1), 1, 1)"
# O2 formula is unchanged: =AVERAGE(OFFSET(Data!B1, 1, 0, 5, 1))
$ws.Range("P2").Value = "Data!C3"
$ws.Range("Q2").Value = "Bob"

# --- Row 3 ---
$ws.Range("N3").Formula = '=INDIRECT("Data!A" & 2)'
# O3 formula is unchanged: =COUNT(INDIRECT("Data!B:B"))
$ws.Range("P3").Value = "Data!A1:C3"
$ws.Range("Q3").Value = $true

# --- Row 4 ---
$ws.Range("N4").Formula = '=INDIRECT("Data!" & CHAR(66) & "2")'
# O4 formula is unchanged: =MAX(INDEX(Data!A1:E6, 0, 4))
$ws.Range("P4").Value = "InvalidSheet!A1"
$ws.Range("Q4").Value = "#REF!"

# --- Row 5 ---
$ws.Range("P5").Value = "'"
$ws.Range("Q5").Value = "#VALUE!"

# --- Row 6 ---
$ws.Range("P6").Value = "Data!A:A"

# --- Row 7 ---
$ws.Range("P7").Value = "Data!1:1"

# The old summary label in row 20 is removed entirely.
$ws.Range("O20").ClearContents()
